$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row ("Density" in A1, "Hardness" in B1) moves one column to the
# right (Density -> B1, Hardness -> C1), leaving A1 blank. The data rows
# (A2:C37) are untouched.
$ws.Range("C1").Value2 = $ws.Range("B1").Value2
$ws.Range("B1").Value2 = $ws.Range("A1").Value2
$ws.Range("A1").Clear()

# C1 used to be an empty, border-only formatted cell; now that it holds the
# "Hardness" header it should carry the same formatting as the other header
# cell (B1) rather than its old corner-border-only style.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Leave the selection on the new first header cell, like the author did.
$ws.Range("B1").Select()
